$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply header style (s="1") to new columns O1:R1 by copying format from N1
$ws.Range("N1").Copy() | Out-Null
$ws.Range("O1:R1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# New header values
$ws.Range("O1").Value = 13
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("R1").Value = 16

# Row 2
$ws.Range("C2").Value = 1.014183356053373
$ws.Range("D2").Value = 1.02960499601579
$ws.Range("E2").Value = 1.026714355959277
$ws.Range("F2").Value = 1.033223717542897
$ws.Range("I2").Value = 1.047549860943408
$ws.Range("J2").Value = 1.035947965810618
$ws.Range("K2").Value = 1.040666736219084
$ws.Range("L2").Value = 1.037813652017308
$ws.Range("M2").Value = 1.044238812832876
$ws.Range("N2").Value = 1.005712725503999
$ws.Range("O2").Value = 1.03
$ws.Range("P2").Value = 1.043584436978209
$ws.Range("Q2").Value = 1.02
$ws.Range("R2").Value = 1.039825185367439

# Row 3
$ws.Range("C3").Value = 1.017618138455724
$ws.Range("D3").Value = 1.031829155068506
$ws.Range("E3").Value = 1.029424912223472
$ws.Range("F3").Value = 1.035943729434249
$ws.Range("I3").Value = 1.048287371232498
$ws.Range("J3").Value = 1.03763977196963
$ws.Range("K3").Value = 1.042070350743421
$ws.Range("L3").Value = 1.039694645519829
$ws.Range("M3").Value = 1.046136466712863
$ws.Range("N3").Value = 1.005712725503983
$ws.Range("O3").Value = 1.03
$ws.Range("P3").Value = 1.045086294855034
$ws.Range("Q3").Value = 1.02
$ws.Range("R3").Value = 1.040814998128236

# Row 4
$ws.Range("C4").Value = 1.019804469617907
$ws.Range("D4").Value = 1.033248759460858
$ws.Range("E4").Value = 1.031156055927204
$ws.Range("F4").Value = 1.037682818475234
$ws.Range("I4").Value = 1.048749337117673
$ws.Range("J4").Value = 1.038714792632806
$ws.Range("K4").Value = 1.042961802171127
$ws.Range("L4").Value = 1.040892452376112
$ws.Range("M4").Value = 1.047346761358246
$ws.Range("N4").Value = 1.005712725503983
$ws.Range("O4").Value = 1.03
$ws.Range("P4").Value = 1.046044155249639
$ws.Range("Q4").Value = 1.02
$ws.Range("R4").Value = 1.041446221883475

# Row 5
$ws.Range("C5").Value = 1.020717999280498
$ws.Range("D5").Value = 1.033844969623259
$ws.Range("E5").Value = 1.031880910081497
$ws.Range("F5").Value = 1.038410949946741
$ws.Range("I5").Value = 1.048941992463206
$ws.Range("J5").Value = 1.039164879953905
$ws.Range("K5").Value = 1.043336576979618
$ws.Range("L5").Value = 1.041393856530856
$ws.Range("M5").Value = 1.047853331149597
$ws.Range("N5").Value = 1.005712725503983
$ws.Range("O5").Value = 1.03
$ws.Range("P5").Value = 1.046445068112934
$ws.Range("Q5").Value = 1.02
$ws.Range("R5").Value = 1.041718371809252

# Row 6
$ws.Range("C6").Value = 1.020874344338575
$ws.Range("D6").Value = 1.033949668675632
$ws.Range("E6").Value = 1.032005228065127
$ws.Range("F6").Value = 1.038535240531031
$ws.Range("I6").Value = 1.048976624712768
$ws.Range("J6").Value = 1.039243550799232
$ws.Range("K6").Value = 1.043404080863895
$ws.Range("L6").Value = 1.041480670173738
$ws.Range("M6").Value = 1.047940439799258
$ws.Range("N6").Value = 1.005712725503983
$ws.Range("O6").Value = 1.03
$ws.Range("P6").Value = 1.04651400820294
$ws.Range("Q6").Value = 1.02
$ws.Range("R6").Value = 1.041774826575486

# Row 7
$ws.Range("C7").Value = 1.019825978490934
$ws.Range("D7").Value = 1.033269774552759
$ws.Range("E7").Value = 1.031173617703224
$ws.Range("F7").Value = 1.037698799463838
$ws.Range("I7").Value = 1.048758625276769
$ws.Range("J7").Value = 1.038729869299564
$ws.Range("K7").Value = 1.042979710922002
$ws.Range("L7").Value = 1.040906937478265
$ws.Range("M7").Value = 1.047359710233715
$ws.Range("N7").Value = 1.005712725503983
$ws.Range("O7").Value = 1.03
$ws.Range("P7").Value = 1.046054403338268
$ws.Range("Q7").Value = 1.02
$ws.Range("R7").Value = 1.041478950552006

# Row 8
$ws.Range("C8").Value = 1.015363368899795
$ws.Range("D8").Value = 1.030376905546599
$ws.Range("E8").Value = 1.027645027268675
$ws.Range("F8").Value = 1.034155168692917
$ws.Range("I8").Value = 1.047810624883572
$ws.Range("J8").Value = 1.036535108735761
$ws.Range("K8").Value = 1.041160542045256
$ws.Range("L8").Value = 1.038463110113291
$ws.Range("M8").Value = 1.044891538639924
$ws.Range("N8").Value = 1.005712725503983
$ws.Range("O8").Value = 1.03
$ws.Range("P8").Value = 1.044101023309334
$ws.Range("Q8").Value = 1.02
$ws.Range("R8").Value = 1.040197208212313

# Row 9
$ws.Range("C9").Value = 1.007196721866025
$ws.Range("D9").Value = 1.025103885255404
$ws.Range("E9").Value = 1.021233951535125
$ws.Range("F9").Value = 1.027734450444004
$ws.Range("I9").Value = 1.046010393896315
$ws.Range("J9").Value = 1.032497880065799
$ws.Range("K9").Value = 1.037802388598828
$ws.Range("L9").Value = 1.033991666804561
$ws.Range("M9").Value = 1.040393029151065
$ws.Range("N9").Value = 1.005712725503983
$ws.Range("O9").Value = 1.03
$ws.Range("P9").Value = 1.040540763620017
$ws.Range("Q9").Value = 1.02
$ws.Range("R9").Value = 1.037819611778609

# Row 10
$ws.Range("C10").Value = 1.00161061544291
$ws.Range("D10").Value = 1.021530047972782
$ws.Range("E10").Value = 1.016922906200749
$ws.Range("F10").Value = 1.02347722824707
$ws.Range("I10").Value = 1.044760355763563
$ws.Range("J10").Value = 1.029760476754885
$ws.Range("K10").Value = 1.035521155361957
$ws.Range("L10").Value = 1.030993190037869
$ws.Range("M10").Value = 1.03743514888731
$ws.Range("N10").Value = 1.005712725503983
$ws.Range("O10").Value = 1.03
$ws.Range("P10").Value = 1.038251248591232
$ws.Range("Q10").Value = 1.02
$ws.Range("R10").Value = 1.036223433439261

# Row 11
$ws.Range("C11").Value = 0.9996788110486238
$ws.Range("D11").Value = 1.020366602136077
$ws.Range("E11").Value = 1.01576727935724
$ws.Range("F11").Value = 1.022743034855378
$ws.Range("I11").Value = 1.044462422849727
$ws.Range("J11").Value = 1.029064957964011
$ws.Range("K11").Value = 1.03491439029354
$ws.Range("L11").Value = 1.03039785012578
$ws.Range("M11").Value = 1.037248428539088
$ws.Range("N11").Value = 1.005712725503983
$ws.Range("O11").Value = 1.03
$ws.Range("P11").Value = 1.03853965749419
$ws.Range("Q11").Value = 1.02
$ws.Range("R11").Value = 1.03582723287609

# Row 12
$ws.Range("C12").Value = 0.9991542631765612
$ws.Range("D12").Value = 1.020078315937481
$ws.Range("E12").Value = 1.015607424133064
$ws.Range("F12").Value = 1.022881631540244
$ws.Range("I12").Value = 1.04444559627778
$ws.Range("J12").Value = 1.028994459282614
$ws.Range("K12").Value = 1.034830698278081
$ws.Range("L12").Value = 1.030441597064749
$ws.Range("M12").Value = 1.037583204119462
$ws.Range("N12").Value = 1.005712725503983
$ws.Range("O12").Value = 1.03
$ws.Range("P12").Value = 1.039130731075317
$ws.Range("Q12").Value = 1.02
$ws.Range("R12").Value = 1.035768060044874

# Row 13
$ws.Range("C13").Value = 0.9996883128196872
$ws.Range("D13").Value = 1.020458014951902
$ws.Range("E13").Value = 1.016213770200057
$ws.Range("F13").Value = 1.02372098389838
$ws.Range("I13").Value = 1.044655314667162
$ws.Range("J13").Value = 1.029413200351472
$ws.Range("K13").Value = 1.035160739313005
$ws.Range("L13").Value = 1.030993781716864
$ws.Range("M13").Value = 1.038364856150775
$ws.Range("N13").Value = 1.005712725503983
$ws.Range("O13").Value = 1.03
$ws.Range("P13").Value = 1.040025281726001
$ws.Range("Q13").Value = 1.02
$ws.Range("R13").Value = 1.035998919776164

# Row 14
$ws.Range("C14").Value = 1.000541118860355
$ws.Range("D14").Value = 1.021027384616354
$ws.Range("E14").Value = 1.016982926649607
$ws.Range("F14").Value = 1.024623110753671
$ws.Range("I14").Value = 1.044904649579157
$ws.Range("J14").Value = 1.029927025879195
$ws.Range("K14").Value = 1.035579822753965
$ws.Range("L14").Value = 1.031608081061101
$ws.Range("M14").Value = 1.039111525853499
$ws.Range("N14").Value = 1.005712725503983
$ws.Range("O14").Value = 1.03
$ws.Range("P14").Value = 1.040788708780328
$ws.Range("Q14").Value = 1.02
$ws.Range("R14").Value = 1.036296647496772

# Row 15
$ws.Range("C15").Value = 1.000984764213074
$ws.Range("D15").Value = 1.021317805398004
$ws.Range("E15").Value = 1.017347094180174
$ws.Range("F15").Value = 1.025013247978456
$ws.Range("I15").Value = 1.0450185770451
$ws.Range("J15").Value = 1.030165791858768
$ws.Range("K15").Value = 1.035778834613503
$ws.Range("L15").Value = 1.031878967267644
$ws.Range("M15").Value = 1.039408983711487
$ws.Range("N15").Value = 1.005712725503983
$ws.Range("O15").Value = 1.03
$ws.Range("P15").Value = 1.041061327429299
$ws.Range("Q15").Value = 1.02
$ws.Range("R15").Value = 1.036443217645081

# Row 16
$ws.Range("C16").Value = 1.003241452106623
$ws.Range("D16").Value = 1.022753867397596
$ws.Range("E16").Value = 1.019055866307007
$ws.Range("F16").Value = 1.026684497267849
$ws.Range("I16").Value = 1.045518486920177
$ws.Range("J16").Value = 1.031250675794345
$ws.Range("K16").Value = 1.036688589261859
$ws.Range("L16").Value = 1.033053685652361
$ws.Range("M16").Value = 1.040552795013197
$ws.Range("N16").Value = 1.005712725503983
$ws.Range("O16").Value = 1.03
$ws.Range("P16").Value = 1.041926709876184
$ws.Range("Q16").Value = 1.02
$ws.Range("R16").Value = 1.037089608080677

# Row 17
$ws.Range("C17").Value = 1.004528881518397
$ws.Range("D17").Value = 1.023560317645524
$ws.Range("E17").Value = 1.019966153044672
$ws.Range("F17").Value = 1.027490489754257
$ws.Range("I17").Value = 1.045770053551003
$ws.Range("J17").Value = 1.03181412712665
$ws.Range("K17").Value = 1.037167486019163
$ws.Range("L17").Value = 1.033632931356782
$ws.Range("M17").Value = 1.041033088992319
$ws.Range("N17").Value = 1.005712725503983
$ws.Range("O17").Value = 1.03
$ws.Range("P17").Value = 1.04217742927995
$ws.Range("Q17").Value = 1.02
$ws.Range("R17").Value = 1.03743079031321

# Row 18
$ws.Range("C18").Value = 1.005087103835607
$ws.Range("D18").Value = 1.023881820282085
$ws.Range("E18").Value = 1.020240455643533
$ws.Range("F18").Value = 1.027567949169852
$ws.Range("I18").Value = 1.04582114786013
$ws.Range("J18").Value = 1.031959294571117
$ws.Range("K18").Value = 1.037300673498878
$ws.Range("L18").Value = 1.033718725038424
$ws.Range("M18").Value = 1.040927206095553
$ws.Range("N18").Value = 1.005712725503983
$ws.Range("O18").Value = 1.03
$ws.Range("P18").Value = 1.041856264780345
$ws.Range("Q18").Value = 1.019999999999999
$ws.Range("R18").Value = 1.03751334031095

# Row 19
$ws.Range("C19").Value = 1.004997547919361
$ws.Range("D19").Value = 1.023783848820141
$ws.Range("E19").Value = 1.019943803485929
$ws.Range("F19").Value = 1.026982493754332
$ws.Range("I19").Value = 1.045703195243612
$ws.Range("J19").Value = 1.031739542817897
$ws.Range("K19").Value = 1.037141881305067
$ws.Range("L19").Value = 1.033364190701106
$ws.Range("M19").Value = 1.040289043124148
$ws.Range("N19").Value = 1.005712725503983
$ws.Range("O19").Value = 1.03
$ws.Range("P19").Value = 1.041026447802649
$ws.Range("Q19").Value = 1.02
$ws.Range("R19").Value = 1.037407454820771

# Row 20
$ws.Range("C20").Value = 1.003082011299808
$ws.Range("D20").Value = 1.022487158336428
$ws.Range("E20").Value = 1.018059680181085
$ws.Range("F20").Value = 1.02459773988723
$ws.Range("I20").Value = 1.045104645411496
$ws.Range("J20").Value = 1.030494641382415
$ws.Range("K20").Value = 1.036145796284305
$ws.Range("L20").Value = 1.031792235066443
$ws.Range("M20").Value = 1.038221430929505
$ws.Range("N20").Value = 1.005712725503983
$ws.Range("O20").Value = 1.03
$ws.Range("P20").Value = 1.038862956177606
$ws.Range("Q20").Value = 1.02
$ws.Range("R20").Value = 1.036707101155211

# Row 21
$ws.Range("C21").Value = 0.9987923555176507
$ws.Range("D21").Value = 1.019736216277605
$ws.Range("E21").Value = 1.014698180174562
$ws.Range("F21").Value = 1.021202757668253
$ws.Range("I21").Value = 1.044103663233695
$ws.Range("J21").Value = 1.028340654766958
$ws.Range("K21").Value = 1.034353084731545
$ws.Range("L21").Value = 1.029406340793303
$ws.Range("M21").Value = 1.035793271427215
$ws.Range("N21").Value = 1.005712725503983
$ws.Range("O21").Value = 1.03
$ws.Range("P21").Value = 1.036900354759852
$ws.Range("Q21").Value = 1.02
$ws.Range("R21").Value = 1.03544281112834

# Row 22
$ws.Range("C22").Value = 0.9960690110119086
$ws.Range("D22").Value = 1.017992302021869
$ws.Range("E22").Value = 1.012591361465419
$ws.Range("F22").Value = 1.019104141690186
$ws.Range("I22").Value = 1.043464639902419
$ws.Range("J22").Value = 1.026982713937904
$ws.Range("K22").Value = 1.033215560233594
$ws.Range("L22").Value = 1.027917440857421
$ws.Range("M22").Value = 1.034306409824952
$ws.Range("N22").Value = 1.005712725503983
$ws.Range("O22").Value = 1.03
$ws.Range("P22").Value = 1.035723596602717
$ws.Range("Q22").Value = 1.02
$ws.Range("R22").Value = 1.03462503731294

# Row 23
$ws.Range("C23").Value = 0.9975087862902652
$ws.Range("D23").Value = 1.018907581522124
$ws.Range("E23").Value = 1.013703751580353
$ws.Range("F23").Value = 1.02021350552518
$ws.Range("I23").Value = 1.043799043744792
$ws.Range("J23").Value = 1.027696742314592
$ws.Range("K23").Value = 1.033809195442633
$ws.Range("L23").Value = 1.028701895678595
$ws.Range("M23").Value = 1.035091096118445
$ws.Range("N23").Value = 1.005712725503983
$ws.Range("O23").Value = 1.03
$ws.Range("P23").Value = 1.036344627071652
$ws.Range("Q23").Value = 1.02
$ws.Range("R23").Value = 1.0350351345549

# Row 24
$ws.Range("C24").Value = 1.003087677411979
$ws.Range("D24").Value = 1.022475357083513
$ws.Range("E24").Value = 1.01803115934868
$ws.Range("F24").Value = 1.024531261256361
$ws.Range("I24").Value = 1.045086345477391
$ws.Range("J24").Value = 1.030467385618602
$ws.Range("K24").Value = 1.036118950744727
$ws.Range("L24").Value = 1.031748857752922
$ws.Range("M24").Value = 1.038140853607403
$ws.Range("N24").Value = 1.005712725503983
$ws.Range("O24").Value = 1.03
$ws.Range("P24").Value = 1.038758315573659
$ws.Range("Q24").Value = 1.02
$ws.Range("R24").Value = 1.036660745687334

# Row 25
$ws.Range("C25").Value = 1.009358530367089
$ws.Range("D25").Value = 1.026506519972981
$ws.Range("E25").Value = 1.022926296109941
$ws.Range("F25").Value = 1.029424946869336
$ws.Range("I25").Value = 1.046500829226422
$ws.Range("J25").Value = 1.033574995758749
$ws.Range("K25").Value = 1.038706833675977
$ws.Range("L25").Value = 1.035178838325978
$ws.Range("M25").Value = 1.04158302692431
$ws.Range("N25").Value = 1.005712725503983
$ws.Range("O25").Value = 1.03
$ws.Range("P25").Value = 1.041482566402756
$ws.Range("Q25").Value = 1.02
$ws.Range("R25").Value = 1.038487649437744
